# Fruta / hortaliza, semanal
# Rotate the weekly price data among rows 3, 4 and 5:
#   new row3 = old row5, new row4 = old row3, new row5 = old row4
# Only columns D, L, M, N, O, P, R, S differ between the rows; the
# remaining columns (A, B, C, E-K, Q, T) are identical across rows 3-5
# and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that rotate, for rows 3-5
$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

$row3 = @{}
$row4 = @{}
$row5 = @{}

foreach ($col in $cols) {
    $row3[$col] = $ws.Range("$col`3").Value2
    $row4[$col] = $ws.Range("$col`4").Value2
    $row5[$col] = $ws.Range("$col`5").Value2
}

# Apply rotation: row3 <- row5, row4 <- row3(old), row5 <- row4(old)
foreach ($col in $cols) {
    $ws.Range("$col`3").Value2 = $row5[$col]
    $ws.Range("$col`4").Value2 = $row3[$col]
    $ws.Range("$col`5").Value2 = $row4[$col]
}
